# Adds a third "Controller" parameter/notification flavour (columns R,S / U,V)
# to the Hilfe sheet, mirroring the existing "Konstruktor Laden Params" (L/M)
# and "Konstruktor Laden" (I/J) columns but targeting a Controller-typed
# parameter (…Controller = xyzController;) for the Events/Notifications work
# between Controller-Data and Controller-HD.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Block 1 header rows (rows 1-6), columns R (params) / U (assignment) ----
$ws.Range("R1").Value = "CharHolder"
$ws.Range("R2").Value = "Konstruktor Laden Params"
$ws.Range("R4").Value = "Controller"
$ws.Range("R5").Formula = '=","'

$ws.Range("S3").Value = " "
$ws.Range("S4").Value = " "
$ws.Range("S5").Value = " "
$ws.Range("S6").Value = " "

$ws.Range("U1").Value = "CharHolder"
$ws.Range("U2").Value = "Konstruktor Laden"
$ws.Range("U4").Formula = '="Controller = "'
$ws.Range("U5").Formula = '="Controller;"'

$ws.Range("V3").Value = " "
$ws.Range("V4").Value = " "
$ws.Range("V5").Value = " "
$ws.Range("V6").Value = " "

# ---- Block 1 data rows (7-16): ObservableCollection<...> entities ----
$block1 = @(
    @{Row=7;  Name="Handlung"},
    @{Row=8;  Name="Fertigkeit"},
    @{Row=9;  Name="Attribut"},
    @{Row=10; Name="Item"},
    @{Row=11; Name="Munition"},
    @{Row=12; Name="Implantat"},
    @{Row=13; Name="Vorteil"},
    @{Row=14; Name="Nachteil"},
    @{Row=15; Name="Connection"},
    @{Row=16; Name="Sin"}
)

foreach ($entry in $block1) {
    $r = $entry.Row
    $ws.Range("R$r").Value = " "
    $ws.Range("S$r").Formula = "=R`$3&`$C$r&R`$4&R`$5"
    $ws.Range("U$r").Value = " "
    $ws.Range("V$r").Formula = "=U`$3&`$C$r&U`$4&LOWER(`$C$r)&U`$5"
}

# ---- Block 1 trailing spacer row (17) ----
$ws.Range("R17").Value = " "
$ws.Range("S17").Value = " "
$ws.Range("U17").Value = " "
$ws.Range("V17").Value = " "

# ---- Block 2 header rows (18-21), same columns, second entity family ----
$ws.Range("S18").Value = " "
$ws.Range("V18").Value = " "

$ws.Range("R19").Value = "Controller"
$ws.Range("S19").Value = " "
$ws.Range("U19").Formula = '="Controller = "'
$ws.Range("V19").Value = " "

$ws.Range("R20").Formula = '=","'
$ws.Range("S20").Value = " "
$ws.Range("U20").Formula = '="Controller;"'
$ws.Range("V20").Value = " "

$ws.Range("R21").Value = " "
$ws.Range("S21").Value = " "
$ws.Range("U21").Value = " "
$ws.Range("V21").Value = " "

# ---- Block 2 data rows (22-27): single-instance entities ----
$block2 = @(
    @{Row=22; Name="Nahkampfwaffe"},
    @{Row=23; Name="Fernkampfwaffe"},
    @{Row=24; Name="Kommlink"},
    @{Row=25; Name="CyberDeck"},
    @{Row=26; Name="Vehikel"},
    @{Row=27; Name="Panzerung"}
)

foreach ($entry in $block2) {
    $r = $entry.Row
    $ws.Range("R$r").Value = " "
    $ws.Range("S$r").Formula = "=R`$18&`$C$r&R`$19&R`$20"
    $ws.Range("U$r").Value = " "
    $ws.Range("V$r").Formula = "=U`$18&`$C$r&U`$19&LOWER(`$C$r)&U`$20"
}

# ---- Selection / view bookkeeping to mirror the authored edit ----
$ws.Range("V22:V27").Select()
